$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "DNET COMMUNICATIONS"
$ws.Range("B4").Value = "'"
$ws.Range("C4").Value = "OS"
$ws.Range("D4").Value = "'9042017010"
$ws.Range("E4").Value = "2025-12-02 14:16"
$ws.Range("F4").Value = "'"
$ws.Range("G4").Value = "2025-12"

# Clear the "quote prefix" style artifact left behind by the leading
# apostrophe (used above to force text storage for numeric-looking /
# empty values) so the new row's cells carry no explicit style, matching
# the unstyled data rows already in the sheet.
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("F4").Style = "Normal"
